# "Generate Report for handoff"
#
# The handoff transform failed for both locales, so this report run no
# longer has a handoff file / handoff datetime / dependency info to show.
#
# "Status" ("Ready for handoff") becomes "Handoff transform failed"
# everywhere it is shown - the Overview sheet's per-locale status cells
# (B2/C2) as well as each locale sheet's own Status cell (B2) - since it
# is literally the same piece of text reused in three places.
#
# For each locale sheet ("zh-cn" and "de-de"):
#   - Status (B2)                   -> "Handoff transform failed"
#   - Latest Handoff File (C2)      -> cleared (no handoff file was produced)
#   - Latest Handoff Datetime (D2)  -> reset to the zero datetime
#   - Latest Handback DateTime (G2) -> reset to the zero datetime
#   - Handoff Reason (H2)           -> "Ignored"
# Row 3 (the .localization-config row) keeps the same displayed values.

$wb = $excel.ActiveWorkbook

$newStatus = "Handoff transform failed"

$mdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/cca2d96b7602f1a8cb69a82268aa4f455d6fbf91/e2e/74e515ba-9fde-4401-933b-8f97ec84b94b.md"
$mdDisplay = "74e515ba-9fde-4401-933b-8f97ec84b94b.md"
$cfgAddress = "https://github.com/OpenLocalizationTest/oltest/blob/cca2d96b7602f1a8cb69a82268aa4f455d6fbf91/.localization-config"
$cfgDisplay = ".localization-config"

# --- Overview sheet: both locale status columns show the new status ------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus

# --- Per-locale detail sheets ---------------------------------------------
$locales = @("zh-cn", "de-de")

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale)

    # --- cell values -------------------------------------------------------
    $ws.Range("B2").Value = $newStatus
    $ws.Range("C2").Clear()
    $ws.Range("D2").Value = "0001-01-01 00:00:00"
    $ws.Range("G2").Value = "0001-01-01 00:00:00"
    $ws.Range("H2").Value = "Ignored"

    $ws.Range("D3").Value = "0001-01-01 00:00:00"
    $ws.Range("G3").Value = "0001-01-01 00:00:00"
    $ws.Range("H3").Value = "Ignored"

    # --- hyperlinks ----------------------------------------------------
    # The host only supports clearing the whole per-sheet Hyperlinks
    # collection at once (no single-item delete), so wipe it and recreate
    # only the two links that survive (A2, A3). This naturally reclaims
    # rId3 for the .localization-config link since the handoff file's
    # hyperlink (old C2 / old rId3) no longer exists.
    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), $mdAddress, "", "", $mdDisplay)
    $ws.Hyperlinks.Add($ws.Range("A3"), $cfgAddress, "", "", $cfgDisplay)
}
